# Update standard cooldown values for all triggers.
# New cooldowns: Basic Kill 10s, Parry 5s, Dismemberment 10s, Critical 10s,
# Decapitation 10s, Last Enemy 20s (was 0), Last Stand 60s (was 45s).
# The allowed-values list (column E on "Menu Mock" and column B on
# "Providers") also changes: "19.6s" -> "20.0s" and a new "60.0s" entry is
# inserted right after "45.0s".

$wb = $excel.ActiveWorkbook

# Original list had "19.6s" (now "20.0s") and was missing "60.0s" (now
# inserted between "45.0s" and "67.5s").
$newList = "0s | 1.6s | 2.0s | 2.4s | 2.5s | 2.8s | 3.0s | 3.5s | 4.0s | 4.2s | 4.9s | 5.0s | 6.0s | 7.0s | 7.2s | 7.5s | 8.0s | 9.0s | 10.0s | 10.5s | 11.2s | 12.6s | 14.0s | 17.5s | 18.0s | 20.0s | 22.5s | 24.5s | 27.0s | 28.0s | 31.5s | 45.0s | 60.0s | 67.5s | 81.0s | 90.0s | 126.0s | 157.5s"

# --- "Menu Mock" sheet: per-trigger standard cooldown (column D) and the
#     shared allowed-values list (column E) on each trigger's Cooldown row.
$ws = $wb.Worksheets.Item("Menu Mock")

$rows = @(45, 54, 63, 72, 81, 89, 97)
$newCooldowns = @{
    45 = "10.0s"  # CategoryCustomBasic         (Basic Kill)
    54 = "10.0s"  # CategoryCustomCritical      (Critical)
    63 = "10.0s"  # CategoryCustomDismemberment (Dismemberment)
    72 = "10.0s"  # CategoryCustomDecapitation  (Decapitation)
    81 = "20.0s"  # CategoryCustomLastEnemy     (Last Enemy, was 0)
    89 = "60.0s"  # CategoryCustomLastStand     (Last Stand, was 45s)
    97 = "5.0s"   # CategoryCustomParry         (Parry)
}

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $newCooldowns[$r]
    $ws.Range("E$r").Value = $newList
}

# --- "Providers" sheet: same allowed-values list mirrored in B10.
$ws2 = $wb.Worksheets.Item("Providers")
$ws2.Range("B10").Value = $newList
